{"js": "// MSM Legal Answer Template (Wills) - fix the answer table's row heights so\n// the PDF export stops cramming the date/contact rows: give the data rows\n// (everything under the \"Question\"/\"Answer\" header row) a fixed (\"exact\")\n// height instead of the default \"at least\" sizing.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// Row 0 is the \"Question\" / \"Answer\" header row and keeps its original\n// height. Rows 1-7 are the data rows being re-sized; row index 2 (\"Does the\n// Main Applicant speak another language...\") is taller because its answer\n// text wraps onto more lines.\nconst newHeightsInTwips = {\n  1: 567, // Legal Issue\n  2: 851, // Does the Main Applicant speak another language...\n  3: 567, // First Name\n  4: 567, // Last Name\n  5: 567, // Date of Birth\n  6: 567, // Best Contact Number\n  7: 567, // Email Address\n};\n\nconst twipsPerPoint = 20;\n\nfor (const [indexStr, heightTwips] of Object.entries(newHeightsInTwips)) {\n  const row = rows.items[Number(indexStr)];\n\n  // The public Word.js TableRow object only exposes `preferredHeight`\n  // (points); it has no public setter for the height *rule*. Drop to the\n  // row's underlying OM setter (the same primitive `preferredHeight`\n  // itself calls internally) to flip the rule to \"exact\", matching\n  // Word's WdRowHeightRule.wdRowHeightExactly used by the desktop app\n  // when a user sets an explicit row height.\n  row._omSet(\"HeightRule\", \"exact\", \"Row\");\n  row.preferredHeight = heightTwips / twipsPerPoint;\n}\n\nawait context.sync();\n", "ps1": "\n# MSM Legal Answer Template (Wills) - fix row heights / date formatting spacing in answer table.\n# The answer table's data rows were given an exact height so the PDF export\n# stops collapsing/overlapping the rows (cantSplit rows need hRule=exact).\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Twips-per-point is 20, Word COM Height/ Rows.Height are expressed in points.\n$wdRowHeightExactly = 2\n\n# Row 2 \"Legal Issue\" -> 567 twips (28.35 pt), exact\n$r = $t.Rows.Item(2)\n$r.HeightRule = $wdRowHeightExactly\n$r.Height = 567 / 20\n\n# Row 3 \"Does the Main Applicant speak another language...\" -> 851 twips (42.55 pt), exact\n$r = $t.Rows.Item(3)\n$r.HeightRule = $wdRowHeightExactly\n$r.Height = 851 / 20\n\n# Row 4 \"First Name\" -> 567 twips, exact\n$r = $t.Rows.Item(4)\n$r.HeightRule = $wdRowHeightExactly\n$r.Height = 567 / 20\n\n# Row 5 \"Last Name\" -> 567 twips, exact\n$r = $t.Rows.Item(5)\n$r.HeightRule = $wdRowHeightExactly\n$r.Height = 567 / 20\n\n# Row 6 \"Date of Birth\" -> 567 twips, exact\n$r = $t.Rows.Item(6)\n$r.HeightRule = $wdRowHeightExactly\n$r.Height = 567 / 20\n\n# Row 7 \"Best Contact Number\" -> 567 twips, exact (was 533)\n$r = $t.Rows.Item(7)\n$r.HeightRule = $wdRowHeightExactly\n$r.Height = 567 / 20\n\n# Row 8 \"Email Address\" -> 567 twips, exact (was 533)\n$r = $t.Rows.Item(8)\n$r.HeightRule = $wdRowHeightExactly\n$r.Height = 567 / 20\n"}
